# Updates the "cryptos" price/volume table (columns D = Price, E = Volume(1h))
# for the 51-row data range, matching the upstream GitHub Actions refresh.
#
# All D/E cells in the sheet are stored as plain text (inline/shared strings),
# even when their contents look numeric (e.g. "601.51", "1.00"). A bare
# Range.Value assignment of a numeric-looking string is auto-coerced to a
# Number by Excel (and would also silently drop significant trailing zeros,
# e.g. "1.00" -> 1, "6.60" -> 6.6 under General format), so for any new D
# value that parses as a plain number we prefix it with a leading apostrophe
# to force text entry, then immediately ClearFormats() on that cell so the
# transient "quote prefix" cell style Excel applies doesn't linger behind
# (keeping the cell's style at the original default/"General" look, with no
# per-cell NumberFormat override). E values (which always contain a percent
# sign and padding spaces) and the handful of D values that already contain
# multiple '.' separators (so Excel can't parse them as numbers anyway) are
# assigned directly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.804.98'
$ws.Range("E2").Value = '  +3.16%  '
$ws.Range("D3").Value = '3.208.52'
$ws.Range("E3").Value = '  +2.12%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = "'601.26"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +3.48%  '
$ws.Range("D6").Value = "'158.16"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +7.46%  '
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").Value = '3.209.18'
$ws.Range("E8").Value = '  +2.18%  '
$ws.Range("E9").Value = '  +5.13%  '
$ws.Range("E10").Value = '  +1.86%  '
$ws.Range("E11").Value = '  -2.41%  '
$ws.Range("E12").Value = '  +3.66%  '
$ws.Range("E13").Value = '  +1.62%  '
$ws.Range("D14").Value = "'39.28"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +5.90%  '
$ws.Range("D15").Value = '3.736.22'
$ws.Range("E15").Value = '  +2.19%  '
$ws.Range("D16").Value = '66.803.09'
$ws.Range("E16").Value = '  +3.03%  '
$ws.Range("E17").Value = '  +4.57%  '
$ws.Range("D18").Value = '3.208.03'
$ws.Range("E18").Value = '  +2.17%  '
$ws.Range("E19").Value = '  +0.89%  '
$ws.Range("D20").Value = "'519.15"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +3.97%  '
$ws.Range("E21").Value = '  +0.16%  '
$ws.Range("E22").Value = '  +4.06%  '
$ws.Range("E23").Value = '  +5.74%  '
$ws.Range("D24").Value = "'15.03"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.30%  '
$ws.Range("E25").Value = '  +1.35%  '
$ws.Range("E26").Value = '  -0.16%  '
$ws.Range("E27").Value = '  +3.11%  '
$ws.Range("D28").Value = "'3.02"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +3.78%  '
$ws.Range("E29").Value = '  +11.34%  '
$ws.Range("D30").Value = "'3.10"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +10.16%  '
$ws.Range("D31").Value = "'7.06"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +10.90%  '
$ws.Range("D32").Value = "'28.28"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +2.62%  '
$ws.Range("E33").Value = '  +1.57%  '
$ws.Range("E34").Value = '  +0.02%  '
$ws.Range("D35").Value = "'6.60"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +2.40%  '
$ws.Range("D36").Value = "'525.75"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +12.02%  '
$ws.Range("D37").Value = "'55.06"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.42%  '
$ws.Range("E38").Value = '  +1.18%  '
$ws.Range("E39").Value = '  +1.96%  '
$ws.Range("E40").Value = '  +9.60%  '
$ws.Range("E41").Value = '  +1.44%  '
$ws.Range("E42").Value = '  +2.19%  '
$ws.Range("E43").Value = '  +9.38%  '
$ws.Range("D44").Value = '0.0₃0692'
$ws.Range("E44").Value = '  +15.73%  '
$ws.Range("E45").Value = '  +3.51%  '
$ws.Range("D46").Value = '2.889.25'
$ws.Range("E46").Value = '  -2.99%  '
$ws.Range("E47").Value = '  +2.33%  '
$ws.Range("D48").Value = "'2.43"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +8.67%  '
$ws.Range("E49").Value = '  +2.90%  '
$ws.Range("E50").Value = '  -0.06%  '
$ws.Range("D51").Value = "'2.68"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +9.99%  '
